$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44232
$ws.Range("J2").Value = 16000

$ws.Range("D3").Value = 44230

$ws.Range("D4").Value = 44187

$ws.Range("D5").Value = 44231
$ws.Range("J5").Value = 12000

$ws.Range("D6").Value = 44209
$ws.Range("K6").Value = 2500
$ws.Range("M6").Value = 2750
$ws.Range("P6").Value = 28

$ws.Range("D7").Value = 44167
$ws.Range("J7").Value = 7000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = 3000
$ws.Range("P7").Value = 30

$ws.Range("D8").Value = 44245
$ws.Range("J8").Value = 9000
$ws.Range("O8").Value = "Región Metropolitana"

$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 5000
$ws.Range("K9").Value = 2500
$ws.Range("L9").Value = 2500
$ws.Range("M9").Value = 2500
$ws.Range("P9").Value = 25

$ws.Range("D10").Value = 44210
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 8800
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 2750
$ws.Range("O10").Value = "Provincia de Chacabuco"
$ws.Range("P10").Value = 28

$ws.Range("D11").Value = 44161
$ws.Range("J11").Value = 7000

$ws.Range("D12").Value = 44214
$ws.Range("K12").Value = 3000
$ws.Range("M12").Value = 3000
$ws.Range("P12").Value = 30

$ws.Range("D13").Value = 44162
$ws.Range("J13").Value = 7000

$ws.Range("D14").Value = 44159
$ws.Range("J14").Value = 7000

$ws.Range("D15").Value = 44188
$ws.Range("J15").Value = 12000

$ws.Range("D16").Value = 44160

$ws.Range("D17").Value = 44189
$ws.Range("J17").Value = 16000

$ws.Range("D19").Value = 44204
$ws.Range("J19").Value = 7000

$ws.Range("D20").Value = 44166

$ws.Range("D21").Value = 44215
$ws.Range("J21").Value = 16000

$ws.Range("D22").Value = 44186
$ws.Range("J22").Value = 10000

$ws.Range("D23").Value = 44168

$ws.Range("D24").Value = 44181
